$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.487.17"
$ws.Range("E2").Value = "  +3.90%  "
$ws.Range("D3").Value = "'1.811.34"
$ws.Range("E3").Value = "  +5.25%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'333.61"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.3806"
$ws.Range("E7").Value = "  +2.45%  "
$ws.Range("D8").Value = "'0.3503"
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("D9").Value = "'49.23"
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("D10").Value = "'1.220"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").Value = "'0.07654"
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'22.10"
$ws.Range("E13").Value = "  +10.28%  "
$ws.Range("D14").Value = "'6.548"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "'1.819.59"
$ws.Range("E15").Value = "  +5.53%  "
$ws.Range("D16").Value = "'7.113"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "'0.00001110"
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("D18").Value = "'0.06691"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").Value = "'85.95"
$ws.Range("E19").Value = "  +4.52%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'17.44"
$ws.Range("E21").Value = "  +5.77%  "
$ws.Range("D22").Value = "'6.502"
$ws.Range("E22").Value = "  +6.54%  "
$ws.Range("D23").Value = "'27.556.26"
$ws.Range("E23").Value = "  +4.28%  "
$ws.Range("D24").Value = "'12.97"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").Value = "'2.441"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "'2.628"
$ws.Range("E26").Value = "  +10.66%  "
$ws.Range("D27").Value = "'21.84"
$ws.Range("E27").Value = "  +12.86%  "
$ws.Range("D28").Value = "'1.455"
$ws.Range("E28").Value = "  +4.69%  "
$ws.Range("D29").Value = "'150.89"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "'2.020.66"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("D31").Value = "'134.29"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "'4.084"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.166"
$ws.Range("E33").Value = "  +3.85%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.08738"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'13.54"
$ws.Range("E35").Value = "  +7.18%  "
$ws.Range("D36").Value = "'1.682"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "'5.534"
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").Value = "'0.6939"
$ws.Range("E38").Value = "  +12.81%  "
$ws.Range("D39").Value = "'9.059"
$ws.Range("E39").Value = "  +7.60%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2233"
$ws.Range("E40").Value = "  +3.98%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.06425"
$ws.Range("E41").Value = "  +3.85%  "
$ws.Range("D42").Value = "'0.02372"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("D43").Value = "'1.296"
$ws.Range("E43").Value = "  +5.87%  "
$ws.Range("D44").Value = "'14.65"
$ws.Range("E44").Value = "  +4.71%  "
$ws.Range("D45").Value = "'0.6491"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'3.875"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "'2.160"
$ws.Range("E48").Value = "  +6.51%  "
$ws.Range("D49").Value = "'131.78"
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("D50").Value = "'0.07281"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").Value = "'80.20"
$ws.Range("E51").Value = "  +4.91%  "
